$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

$ws.Range("D5").Value = 2.363070468628005
$ws.Range("D6").Value = 0.07188899101417243
$ws.Range("D7").Value = -0.3481388986083073
$ws.Range("D8").Value = 0.2262012405916103
$ws.Range("D9").Value = 2.492867487854896
$ws.Range("D10").Value = 0.2671958650971716
$ws.Range("D11").Value = 2.418529064907384
$ws.Range("D12").Value = 0.3311118621155796
$ws.Range("D13").Value = 0.4266508955728286
$ws.Range("D14").Value = 0.2009665427509293
$ws.Range("D15").Value = 0.2336306198937148
$ws.Range("D16").Value = 0.2780790558736198
$ws.Range("D17").Value = -0.1564787544150041
$ws.Range("D18").Value = 0.005974167488294789
$ws.Range("D19").Value = 0.4117376532571767
$ws.Range("D20").Value = 0.6240486996118719
$ws.Range("D21").Value = 0.2237575135350566
$ws.Range("D22").Value = 0.03675178589910508
$ws.Range("D23").Value = 0.2366326647980964
$ws.Range("D24").Value = 0.3266854758224005
